$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Session 5 (Dynamic Prog.) mark for the student
$ws.Range("F4").Value = 10

# Teacher comment for the Session 5 mark
$ws.Range("F5").Value = "Perfect, keep on!"
